$d = $word.ActiveDocument

$newText = "A way that I can figure out the solution would be to count from 1 to 1000 in the same manner that she did, but that would be time consuming. "

# Locate the paragraph that ends the "Predicting Fingers" section by matching its
# text (robust against exact-whitespace/index drift).
$idx = 0
$anchorIdx = -1
foreach ($p in $d.Paragraphs) {
    $idx++
    if ($p.Range.Text -like "*same finger every time*") {
        $anchorIdx = $idx
    }
}

if ($anchorIdx -eq -1) {
    throw "Could not locate anchor paragraph"
}

$anchorPara = $d.Paragraphs.Item($anchorIdx)

# Insert a brand-new paragraph right after the anchor, inheriting its
# paragraph/run formatting (Helvetica 13pt, spacing after 260).
$anchorPara.Range.InsertParagraphAfter()

# The freshly-created paragraph now sits at $anchorIdx + 1; fill in its text.
$newPara = $d.Paragraphs.Item($anchorIdx + 1)
$newPara.Range.Text = $newText
